# issue #5: stock data from json to db
#
# The "股票" (stock) worksheet gains a new "category" column (inserted
# between property_category and date) plus two new trailing columns,
# source_file and index, appended after the existing legislator_id
# column. The sheet's used range grows from A1:K3 to A1:N3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)   # 股票

# Insert a blank column before the existing "date" column (I) so that
# date / legislator_name / legislator_id shift right to J / K / L,
# carrying their values, types and formatting with them intact.
$ws.Range("I1").EntireColumn.Insert()

# --- new "category" column (now column I) --------------------------------
$ws.Range("I1").Value = "category"
$ws.Range("I2").Value = "normal"
$ws.Range("I3").Value = "normal"

# --- new trailing columns: source_file (M) and index (N) -----------------
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

$ws.Range("M2").Value = "tmped121"
$ws.Range("N2").Value = 69

$ws.Range("M3").Value = "tmped121"
$ws.Range("N3").Value = 70

# Copy the existing header formatting (bold, centered, bordered) onto the
# two new trailing header cells so they match the rest of row 1.
$ws.Range("L1").Copy()
$ws.Range("M1:N1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
